# Daily attendance processing - 2025-12-16 10:32:48
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Recorded By" column (G): swap "System, <email>" -> "<email>, System"
#     for every recorded session row in the sheet.
$recordedByRows = @(2,3,4,24,25,26,46,47,48,68,69,70,90,91,92,112,113,114,134,135,136,156,157,158,178,179,180,200,201,202,222,223,224,244,245,246)
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# --- 2. Summary block: Missing / Pending session totals shifted by the
#     six sessions that flipped from "Pending" to "Not Recorded" below.
$ws.Range("L7").Value = 24   # Missing Sessions
$ws.Range("L8").Value = 186  # Pending Sessions

# --- 3. Group statistics table: Missing (P) +1 / Pending (Q) -1 for the
#     six groups whose 16/12/2025 session is now overdue ("Not Recorded").
$groupStatRows = @(16,17,18,24,25,26)
foreach ($r in $groupStatRows) {
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 16
}

# --- 4. The six 16/12/2025 sessions themselves move from "Pending"
#     (yellow) to "Not Recorded" (pink/red) now that their date has passed.
#     Re-use an existing "Not Recorded" row's formatting (e.g. row 18) so
#     the pasted style matches the workbook's existing pink style exactly.
$formatSourceRow = 18
$statusRows = @(27,49,71,203,225,247)
foreach ($r in $statusRows) {
    $ws.Range("A$formatSourceRow`:I$formatSourceRow").Copy() | Out-Null
    $ws.Range("A$r`:I$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("I$r").Value = "Not Recorded"
}

$excel.CutCopyMode = 0
